$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at Y (25) -- everything from Y onward shifts right by one.
$ws.Columns("Y:Y").Insert()

# New header cell Y1: "Забрал оригинал", centered horizontally.
$ws.Range("Y1").Value = "Забрал оригинал"
$ws.Range("Y1").HorizontalAlignment = -4108

# New field-code cell Y2: "[application.is_return_original_epk]"
$ws.Range("Y2").Value = "[application.is_return_original_epk]"

# Try to match the original column's width as closely as possible.
$ws.Columns("Y:Y").ColumnWidth = $ws.Columns("X:X").ColumnWidth

# Update the active selection to match the target view.
$ws.Range("X7").Select()
